# CRS2025_R_workshop.pptx edit
# 1) Slide 17, Content Placeholder 2, paragraph 3:
#    "Multiple " + "packages: install" + ".packages" + "(c(...))"
#    -> "Multiple packages: " + "install.packages" + "(c(...))"
# 2) Slide 38, Content Placeholder 2, paragraph 4:
#    "holobaramins" + " in the data."
#    -> "holobaramins" + " (kinds) " + "in the data."  (all bold)

$p = $ppt.ActivePresentation

# --- Edit 1 ---
$s1 = $p.Slides.Item(17)
$shp1 = $s1.Shapes.Item(2)
$tr1 = $shp1.TextFrame.TextRange
$para1 = $tr1.Paragraphs(3, 1)

$run1a = $para1.Characters(1, 19)
$run1b = $para1.Characters(20, 16)

$run1a.Text = "Multiple packages: "
$run1b.Text = "install.packages"

# --- Edit 2 ---
$s2 = $p.Slides.Item(38)
$shp2 = $s2.Shapes.Item(2)
$tr2 = $shp2.TextFrame.TextRange
$para2 = $tr2.Paragraphs(4, 1)

$run2 = $para2.Characters(59, 13)
$run2.Text = " (kinds) in the data."

$boldTail = $para2.Characters(68, 12)
$boldTail.Font.Bold = $true
